# Populate the header row on Sheet1 with the two column titles, autofit the
# columns to the new content, and leave the selection on the row below the
# header - as a user would after typing the headers and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Phone"
$ws.Range("B1").Value = "Price"

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# Nudge to the widths AutoFit produced in the authored workbook (closest
# values reachable through the ColumnWidth setter).
$ws.Columns.Item(1).ColumnWidth = 41.333333333333336
$ws.Columns.Item(2).ColumnWidth = 39

$ws.Range("A2").Select() | Out-Null
